$d = $word.ActiveDocument

# Helper: append a new paragraph at the very end of the document body, with the
# given text (pass empty string for a blank paragraph). Because the current
# last paragraph ("DROP COLUMN usuariopedido;") carries underline-only
# character formatting (no color/size override) on its paragraph mark, each
# newly inserted paragraph mark - and any text subsequently typed into it -
# naturally inherits that same underline formatting, matching the rest of
# this section of the document.
function Add-Line([string]$text) {
    $r = $d.Content
    $r.Collapse(0)  # wdCollapseEnd
    $r.InsertParagraphAfter()
    if ($text.Length -gt 0) {
        $r = $d.Content
        $r.Collapse(0)  # wdCollapseEnd
        $r.Text = $text
    }
}

Add-Line ""
Add-Line "ALTER TABLE bebidas ADD COLUMN imagem BYTEA;"
Add-Line ""
Add-Line "ALTER TABLe bebidas"
Add-Line "DROP COLUMN qntbebida;"
